$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - values updated
$ws.Range("B3").Value = 0.997793535912601
$ws.Range("C3").Value = 0.9978715767131624
$ws.Range("D3").Value = 0.985285803019826

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor, values updated
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9978428846270629
$ws.Range("C4").Value = 0.9976987213718488
$ws.Range("D4").Value = 0.9742013305079448

# Row 5: AdaBoostRegressor -> MLPRegressor, values updated
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9985109586975921
$ws.Range("C5").Value = 0.9982210298571893
$ws.Range("D5").Value = 0.9979573974830886
